# Generate Report for Handoff
#
# A new source file (21ed720f-7e9e-4c8d-bada-a09f64d3d4ac) became ready for
# handoff. It is inserted as a new row ABOVE the existing
# 63f49407-adcd-4efb-ace2-c3cf4b81369b row on every sheet (Overview, zh-cn,
# de-de) - i.e. the existing row 2 becomes row 3, and the new file's data is
# written into row 2.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"  (File Name | zh-cn | de-de | Latest Handoff Date)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$oldMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/5d207da75a25e77f443371af7a8ff8d692e2d34e/e2e/63f49407-adcd-4efb-ace2-c3cf4b81369b.md"
$newMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/5d207da75a25e77f443371af7a8ff8d692e2d34e/e2e/21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.md"

# Duplicate row 2 down into row 3 - this preserves cell styles (hyperlink
# font etc.) exactly, the values get overwritten/re-linked below.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(2).Insert(-4121)

# The old hyperlink object is still anchored on A2 after the shift; drop it,
# then re-create it (and the one for the new row) in the right order so the
# relationship ids line up: A2 -> rId2, A3 -> rId3.
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A3").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", "21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.md") | Out-Null
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-22-11 08:22:00"

$ws.Hyperlinks.Add($ws.Range("A3"), $oldMdUrl, "", "", "63f49407-adcd-4efb-ace2-c3cf4b81369b.md") | Out-Null
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-21-11 08:21:33"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$oldMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/5d207da75a25e77f443371af7a8ff8d692e2d34e/e2e/63f49407-adcd-4efb-ace2-c3cf4b81369b.md"
$newMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/5d207da75a25e77f443371af7a8ff8d692e2d34e/e2e/21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.md"
$oldXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/52cb1c90da76ff505e8754f72c23883e7ee19648/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/63f49407-adcd-4efb-ace2-c3cf4b81369b.51a3d6c5a4cd58e4fef9e873ef7c2ef0a2201152.zh-cn.xlf"
$newXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/52cb1c90da76ff505e8754f72c23883e7ee19648/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.7b32db215b5030fc8eef5443a292995986e0f93c.zh-cn.xlf"

$ws.Rows.Item(2).Copy()
$ws.Rows.Item(2).Insert(-4121)

$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("D2").Hyperlinks.Delete()
$ws.Range("A3").Hyperlinks.Delete()
$ws.Range("B3").Hyperlinks.Delete()
$ws.Range("D3").Hyperlinks.Delete()

# Row 2 -> new file (21ed720f...)
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", "21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), $newMdUrl, "", "", ".md") | Out-Null
$ws.Range("C2").Value = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("D2"), $newXlfUrl, "", "", "21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.7b32db215b5030fc8eef5443a292995986e0f93c.zh-cn.xlf") | Out-Null
$ws.Range("E2").Value = "2016-03-11 08:21:57"
$ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

# Row 3 -> old file (63f49407...)
$ws.Hyperlinks.Add($ws.Range("A3"), $oldMdUrl, "", "", "63f49407-adcd-4efb-ace2-c3cf4b81369b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $oldMdUrl, "", "", ".md") | Out-Null
$ws.Range("C3").Value = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("D3"), $oldXlfUrl, "", "", "63f49407-adcd-4efb-ace2-c3cf4b81369b.51a3d6c5a4cd58e4fef9e873ef7c2ef0a2201152.zh-cn.xlf") | Out-Null
$ws.Range("E3").Value = "2016-03-11 08:21:30"
$ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$oldMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/5d207da75a25e77f443371af7a8ff8d692e2d34e/e2e/63f49407-adcd-4efb-ace2-c3cf4b81369b.md"
$newMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/5d207da75a25e77f443371af7a8ff8d692e2d34e/e2e/21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.md"
$oldXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01564b4e459b5b7b670ae9fb70e9284aca727d1a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/63f49407-adcd-4efb-ace2-c3cf4b81369b.51a3d6c5a4cd58e4fef9e873ef7c2ef0a2201152.de-de.xlf"
$newXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01564b4e459b5b7b670ae9fb70e9284aca727d1a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.7b32db215b5030fc8eef5443a292995986e0f93c.de-de.xlf"

$ws.Rows.Item(2).Copy()
$ws.Rows.Item(2).Insert(-4121)

$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("D2").Hyperlinks.Delete()
$ws.Range("A3").Hyperlinks.Delete()
$ws.Range("B3").Hyperlinks.Delete()
$ws.Range("D3").Hyperlinks.Delete()

# Row 2 -> new file (21ed720f...)
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", "21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), $newMdUrl, "", "", ".md") | Out-Null
$ws.Range("C2").Value = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("D2"), $newXlfUrl, "", "", "21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.7b32db215b5030fc8eef5443a292995986e0f93c.de-de.xlf") | Out-Null
$ws.Range("E2").Value = "2016-03-11 08:22:00"
$ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

# Row 3 -> old file (63f49407...)
$ws.Hyperlinks.Add($ws.Range("A3"), $oldMdUrl, "", "", "63f49407-adcd-4efb-ace2-c3cf4b81369b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), $oldMdUrl, "", "", ".md") | Out-Null
$ws.Range("C3").Value = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("D3"), $oldXlfUrl, "", "", "63f49407-adcd-4efb-ace2-c3cf4b81369b.51a3d6c5a4cd58e4fef9e873ef7c2ef0a2201152.de-de.xlf") | Out-Null
$ws.Range("E3").Value = "2016-03-11 08:21:33"
$ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

Write-Host "Report generated for handoff."
